$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.494.00"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.648.48"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3798"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.223"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.416"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.414"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001203"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "1.656.15"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06984"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.768"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "23.514.46"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.486"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.908"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.79%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "1.838.01"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.934"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  +4.71%  "
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.030"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02732"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08738"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2455"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.969"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06859"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6913"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.323"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6447"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.271"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07803"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.171"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.46%  "
